# react-redux-async.pptx — "Changed from comments in last teach"
#
# Re-merges several runs that had gotten split apart (typically by
# PowerPoint's spell-check "err" splitting or by a previous partial edit)
# back into single runs with the correct combined text, and bolds/colors
# the word "all" in one bullet.
#
# Helper: replace the text spanning [absStart, absStart+len-1] (1-based,
# whole TextFrame.TextRange-relative) with newText. Because the run
# formatting of the FIRST character in the span is kept by PowerPoint when
# a multi-run selection is retyped, this naturally collapses adjacent runs
# that already share that formatting into one run.
# NB: always pass pre-computed plain variables (not inline parenthesized
# expressions) as the positional arguments here.
function Set-Span {
    param($TextRange, [int]$AbsStart, [int]$Len, [string]$NewText)
    $span = $TextRange.Characters($AbsStart, $Len)
    $span.Text = $NewText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: "Signature is ( p1, p2 ) => null" — merge "( " + "p1, p2 ) =>"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$para3 = $tr3.Paragraphs(6)
$idx3 = $para3.Text.IndexOf("( p1")
$abs3 = $para3.Start + $idx3
$text3 = "( p1, p2 ) =>"
$len3 = $text3.Length
Set-Span $tr3 $abs3 $len3 $text3

# ---------------------------------------------------------------------
# Slide 4: "Now, we put all the commands (whoops.. Actions) in a file"
# Split "all" out into its own bold-red run.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange
$para4 = $tr4.Paragraphs(2)
$idx4 = $para4.Text.IndexOf("all")
$abs4 = $para4.Start + $idx4
$len4 = "all".Length
$wordSpan4 = $tr4.Characters($abs4, $len4)
$wordSpan4.Font.Color.RGB = 255

# ---------------------------------------------------------------------
# Slide 7: three run re-merges in the "Text Placeholder 1" body
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(1).TextFrame.TextRange

# "The saveQuote(  name, quote )"
$p7a = $tr7.Paragraphs(8)
$i7a = $p7a.Text.IndexOf("( ")
$abs7a = $p7a.Start + $i7a
$text7a = "(  name, quote )"
$len7a = $text7a.Length
Set-Span $tr7 $abs7a $len7a $text7a

# "The requestQuote( ) which sets a property ..."
$p7b = $tr7.Paragraphs(10)
$i7b = $p7b.Text.IndexOf(" which sets")
$abs7b = $p7b.Start + $i7b
$text7b = " which sets a property "
$len7b = $text7b.Length
Set-Span $tr7 $abs7b $len7b $text7b

# "The receiveQuote(  name, quote ) which does the old save"
$p7c = $tr7.Paragraphs(12)
$i7c = $p7c.Text.IndexOf("( ")
$abs7c = $p7c.Start + $i7c
$text7c = "(  name, quote )"
$len7c = $text7c.Length
Set-Span $tr7 $abs7c $len7c $text7c

# ---------------------------------------------------------------------
# Slide 8: "QuoteA" + "ctions.js" -> "QuoteActions.js"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(1).TextFrame.TextRange
$para8 = $tr8.Paragraphs(1)
$idx8 = $para8.Text.IndexOf("QuoteA")
$abs8 = $para8.Start + $idx8
$text8 = "QuoteActions.js"
$len8 = $text8.Length
Set-Span $tr8 $abs8 $len8 $text8

# ---------------------------------------------------------------------
# Slide 9: same filename fix, plus "action " + "name  !!!!!" merge
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(1).TextFrame.TextRange
$para9 = $tr9.Paragraphs(1)
$idx9 = $para9.Text.IndexOf("QuoteA")
$abs9 = $para9.Start + $idx9
$text9 = "QuoteActions.js"
$len9 = $text9.Length
Set-Span $tr9 $abs9 $len9 $text9

$grp9 = $s9.Shapes.Item(4)
$tb9 = $grp9.GroupItems.Item(2)
$tr9b = $tb9.TextFrame.TextRange
$idx9b = $tr9b.Text.IndexOf("action ")
$abs9b = $idx9b + 1
$text9b = "action name  !!!!!"
$len9b = $text9b.Length
Set-Span $tr9b $abs9b $len9b $text9b

# ---------------------------------------------------------------------
# Slide 10: "QuoteR" + "educer.js" -> "QuoteReducer.js"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(1).TextFrame.TextRange
$para10 = $tr10.Paragraphs(1)
$idx10 = $para10.Text.IndexOf("QuoteR")
$abs10 = $para10.Start + $idx10
$text10 = "QuoteReducer.js"
$len10 = $text10.Length
Set-Span $tr10 $abs10 $len10 $text10
